$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 166, pushing existing rows 166:262 down to 167:263
$ws.Rows("166:166").Insert()

# Populate the newly inserted row with the new data record
$ws.Range("A166").Value = 10
$ws.Range("B166").Value = "Vega Modelo de Temuco"
$ws.Range("C166").Value = "La Araucanía"
$ws.Range("D166").Value = 45029
$ws.Range("E166").Value = 9
$ws.Range("F166").Value = 100112012
$ws.Range("G166").Value = "Espinaca"
$ws.Range("H166").Value = "Sin especificar"
$ws.Range("I166").Value = "Primera"
$ws.Range("J166").Value = 40
$ws.Range("K166").Value = 12000
$ws.Range("L166").Value = 12000
$ws.Range("M166").Value = 12000
$ws.Range("N166").Value = "$/docena de atados"
$ws.Range("O166").Value = "Región de La Araucanía"
$ws.Range("P166").Value = 4000
$ws.Range("Q166").Value = 3
$ws.Range("R166").Value = "Hortaliza"
